$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header C1: siteid -> sitenumber
$ws.Range("C1").Value = "sitenumber"

# Update lat/lon values for rows 2-11
$ws.Range("A2").Value = 43.21166
$ws.Range("B2").Value = -76.28693800000001

$ws.Range("A3").Value = 32.75469
$ws.Range("B3").Value = -96.77826

$ws.Range("A4").Value = 39.118211
$ws.Range("B4").Value = -94.93859

$ws.Range("A5").Value = 39.736799
$ws.Range("B5").Value = -105.07198

$ws.Range("A6").Value = 30.45373
$ws.Range("B6").Value = -89.01851000000001

$ws.Range("A7").Value = 44.07253
$ws.Range("B7").Value = -123.06642

$ws.Range("A8").Value = 35.704882
$ws.Range("B8").Value = -81.31295900000001

$ws.Range("A9").Value = 18.05673
$ws.Range("B9").Value = -66.72190999999999

$ws.Range("A10").Value = 41.953372
$ws.Range("B10").Value = -87.68844

$ws.Range("A11").Value = 39.735207
$ws.Range("B11").Value = -86.140005
